# Rename the ACC.* sub-statement sheets (and their hidden "_" list-source
# siblings) to be prefixed with "ACC." so they no longer collide with the
# generic REF.* naming, then fix up every formula / data-validation source
# that referenced the old names.

$wb = $excel.ActiveWorkbook

# Old name -> New name, in workbook tab order (sheetId 49 .. 78)
$renames = [ordered]@{
    "ASST.SFT"        = "ACC.ASST_SFT"
    "ASST.SFT_"       = "ACC.ASST_SFT_"
    "ASST.DBT"        = "ACC.ASST_DBT"
    "ASST.DBT_"       = "ACC.ASST_DBT_"
    "ASST.LN"         = "ACC.ASST_LN"
    "ASST.LN_"        = "ACC.ASST_LN_"
    "SHR"             = "ACC.SHR"
    "SHR_"            = "ACC.SHR_"
    "EDR"             = "ACC.EDR"
    "EDR_"            = "ACC.EDR_"
    "ODR"             = "ACC.ODR"
    "ODR_"            = "ACC.ODR_"
    "ASST.NN_FNNCL"   = "ACC.ASST_NN_FNNCL"
    "ASST.NN_FNNCL_"  = "ACC.ASST_NN_FNNCL_"
    "ASST.RMNNG"      = "ACC.ASST_RMNNG"
    "ASST.RMNNG_"     = "ACC.ASST_RMNNG_"
    "LBLTY.SFT"       = "ACC.LBLTY_SFT"
    "LBLTY.SFT_"      = "ACC.LBLTY_SFT_"
    "LBLTY.DBT"       = "ACC.LBLTY_DBT"
    "LBLTY.DBT_"      = "ACC.LBLTY_DBT_"
    "LBLTY.LN"        = "ACC.LBLTY_LN"
    "LBLTY.LN_"       = "ACC.LBLTY_LN_"
    "LBLTY.RMNNG"     = "ACC.LBLTY_RMNNG"
    "LBLTY.RMNNG_"    = "ACC.LBLTY_RMNNG_"
    "HLDR"            = "ACC.HLDR"
    "HLDR_"           = "ACC.HLDR_"
    "RVN"             = "ACC.RVN"
    "RVN_"            = "ACC.RVN_"
    "EXPNS"           = "ACC.EXPNS"
    "EXPNS_"          = "ACC.EXPNS_"
}

foreach ($old in $renames.Keys) {
    $ws = $wb.Worksheets.Item($old)
    $ws.Name = $renames[$old]
}

# The CONTENTS sheet (row 25-39, column A) has HYPERLINK() formulas whose
# string arguments embed the old sheet names literally - update those too.
$contents = $wb.Worksheets.Item("CONTENTS")
$hyperlinkRows = @{
    25 = "ACC.ASST_SFT"
    26 = "ACC.ASST_DBT"
    27 = "ACC.ASST_LN"
    28 = "ACC.SHR"
    29 = "ACC.EDR"
    30 = "ACC.ODR"
    31 = "ACC.ASST_NN_FNNCL"
    32 = "ACC.ASST_RMNNG"
    33 = "ACC.LBLTY_SFT"
    34 = "ACC.LBLTY_DBT"
    35 = "ACC.LBLTY_LN"
    36 = "ACC.LBLTY_RMNNG"
    37 = "ACC.HLDR"
    38 = "ACC.RVN"
    39 = "ACC.EXPNS"
}
foreach ($row in $hyperlinkRows.Keys) {
    $name = $hyperlinkRows[$row]
    $cell = $contents.Range("A$row")
    $cell.Formula = '=HYPERLINK("#' + $name + '!A1", "' + $name + '")'
}

# Each renamed "visible" sheet's data-validation dropdown (on the statement
# entry cells) sources its list from the matching hidden "_" sheet by name -
# update those formula1 references to the new sheet names.
$validationFixups = @(
    @{ Sheet = "ACC.ASST_SFT";        Ranges = @("C4:C5") },
    @{ Sheet = "ACC.ASST_DBT";        Ranges = @("D4:D5") },
    @{ Sheet = "ACC.ASST_LN";         Ranges = @("C4:C5") },
    @{ Sheet = "ACC.SHR";             Ranges = @("D4:D5") },
    @{ Sheet = "ACC.EDR";             Ranges = @("D4:D5") },
    @{ Sheet = "ACC.ODR";             Ranges = @("C4:C5") },
    @{ Sheet = "ACC.ASST_NN_FNNCL";   Ranges = @("A4:A5", "B4:B5", "E4:E5") },
    @{ Sheet = "ACC.ASST_RMNNG";      Ranges = @("A4:A5", "D4:D5") },
    @{ Sheet = "ACC.LBLTY_SFT";       Ranges = @("C4:C5") },
    @{ Sheet = "ACC.LBLTY_DBT";       Ranges = @("D4:D5") },
    @{ Sheet = "ACC.LBLTY_LN";        Ranges = @("C4:C5") },
    @{ Sheet = "ACC.LBLTY_RMNNG";     Ranges = @("A4:A5", "D4:D5") },
    @{ Sheet = "ACC.HLDR";            Ranges = @("D4:D5") },
    @{ Sheet = "ACC.RVN";             Ranges = @("B4:B5", "D4:D5") },
    @{ Sheet = "ACC.EXPNS";           Ranges = @("B4:B5", "D4:D5") }
)

foreach ($fixup in $validationFixups) {
    $ws = $wb.Worksheets.Item($fixup.Sheet)
    foreach ($addr in $fixup.Ranges) {
        $rng = $ws.Range($addr)
        $dv = $rng.Validation
        $formula1 = $dv.Formula1
        if ($formula1 -match "^'([^']+)'!(.*)$") {
            $oldSheetRef = $matches[1]
            $rest = $matches[2]
            if ($renames.Contains($oldSheetRef)) {
                $newSheetRef = $renames[$oldSheetRef]
                $dv.Formula1 = "'" + $newSheetRef + "'!" + $rest
            }
        }
    }
}
